$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NZREG code rows (Code, Description) to append below the existing
# table; the Definition column (C) stays blank, matching the existing
# blank Definition cells used throughout the sheet.
$newRows = @(
    @("94", "Don't Know"),
    @("95", "Refused to Answer"),
    @("96", "Repeated Value"),
    @("97", "Response Unidentifiable"),
    @("98", "Response Outside Scope"),
    @("99", "Not Stated")
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A holds zero-padded numeric codes; prefix with a text quote
    # marker so "94" etc. are stored as text (matching the existing
    # "00".."23" codes) rather than being coerced to numbers, then strip
    # the quote-prefix formatting it introduces.
    $ws.Range("A$r").Value = "'" + $rowData[0]
    $ws.Range("A$r").ClearFormats()

    $ws.Range("B$r").Value = $rowData[1]

    # Column C (Definition) is left blank for these rows, same as all
    # other rows in the table. A bare "" assignment clears/omits the
    # cell entirely, so use the text quote-prefix trick to force an
    # actual empty-text cell, then strip the quote-prefix formatting it
    # introduces.
    $ws.Range("C$r").Value = "'"
    $ws.Range("C$r").ClearFormats()
}
